$wb = $excel.ActiveWorkbook

# 展览 (Worksheets(1))
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 64
$ws.Range("F3").Value = 1032
$ws.Range("F4").Value = 9200
$ws.Range("F5").Value = 189
$ws.Range("F6").Value = 63
$ws.Range("F7").Value = 1955
$ws.Range("F8").Value = 6367
$ws.Range("F9").Value = 616
$ws.Range("F10").Value = 64
$ws.Range("F11").Value = 9409
$ws.Range("F12").Value = 10977
$ws.Range("F13").Value = 1225
$ws.Range("F14").Value = 1131
$ws.Range("F15").Value = 4900
$ws.Range("F17").Value = 441
$ws.Range("F19").Value = 325
$ws.Range("F22").Value = 233
$ws.Range("F23").Value = 869
$ws.Range("F24").Value = 1219
$ws.Range("F25").Value = 851
$ws.Range("F26").Value = 1
$ws.Range("F27").Value = 2016
$ws.Range("F28").Value = 419
$ws.Range("F29").Value = 609
$ws.Range("F30").Value = 2637
$ws.Range("F32").Value = 181
$ws.Range("F33").Value = 1720
$ws.Range("F36").Value = 438
$ws.Range("F37").Value = 43
$ws.Range("F38").Value = 912
$ws.Range("F39").Value = 584
$ws.Range("F40").Value = 13
$ws.Range("F41").Value = 3292
$ws.Range("F44").Value = 502
$ws.Range("F47").Value = 894
$ws.Range("F49").Value = 4199

# 演出 (Worksheets(2))
$ws = $wb.Worksheets.Item(2)
$ws.Range("F11").Value = 65
$ws.Range("F23").Value = 67

# 本地生活 (Worksheets(3))
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 5860

# 全部类型 (Worksheets(4))
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1032
$ws.Range("F3").Value = 9200
$ws.Range("F4").Value = 63
$ws.Range("F7").Value = 6367
$ws.Range("F8").Value = 616
$ws.Range("F9").Value = 9409
$ws.Range("F10").Value = 10977
$ws.Range("F12").Value = 1225
$ws.Range("F13").Value = 1131
$ws.Range("F14").Value = 4900
$ws.Range("F16").Value = 441
$ws.Range("F19").Value = 65
$ws.Range("F21").Value = 233
$ws.Range("F22").Value = 869
$ws.Range("F23").Value = 1219
$ws.Range("F24").Value = 851
$ws.Range("F26").Value = 2016
$ws.Range("F27").Value = 419
$ws.Range("F28").Value = 609
$ws.Range("F29").Value = 2637
$ws.Range("F30").Value = 181
$ws.Range("F31").Value = 1720
$ws.Range("F33").Value = 438
$ws.Range("F36").Value = 43
$ws.Range("F37").Value = 912
$ws.Range("F39").Value = 67
$ws.Range("F40").Value = 584
$ws.Range("F44").Value = 502
$ws.Range("F47").Value = 894
